$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: replace "Example 2" placeholder with real maze-generation log entry
$ws.Range("A5").Value = "Setting up structure for maze generation"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "2/6/2022"
$ws.Range("D5").Value = "Creating UML diagram to plan out the handling of a perfect maze generator"
$ws.Range("D5").WrapText = $true

# Row 6: replace "Example 3" placeholder with real maze-generation log entry
$ws.Range("A6").Value = "Starting Creation of Maze algorithm"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "2/6/2022"
$ws.Range("D6").Value = "Original maze structure creation as setup for Prim's algorithm"
$ws.Range("E6").Value = ""

# Move the active selection to D11, matching where the author was working
$ws.Range("D11").Select()
